$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking strings (e.g. "26.891.20")
# are stored verbatim instead of being parsed into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.891.20"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "1.546.52"
$ws.Range("E3").Value = "  -0.93%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "205.86"
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("D9").Value = "21.30"
$ws.Range("E9").Value = "  -2.13%  "

$ws.Range("E10").Value = "  -0.35%  "

$ws.Range("D11").Value = "0.0854"
$ws.Range("E11").Value = "  -1.25%  "

$ws.Range("D12").Value = "1.766.27"
$ws.Range("E12").Value = "  -1.04%  "

$ws.Range("D13").Value = "1.543.44"
$ws.Range("E13").Value = "  -1.25%  "

$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").Value = "0.511"
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("D16").Value = "26.864.29"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").Value = "61.51"
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("D18").Value = "213.55"
$ws.Range("E18").Value = "  -0.54%  "

$ws.Range("E19").Value = "  +0.27%  "

$ws.Range("E20").Value = "  -1.99%  "

$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("E22").Value = "  -2.22%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("E24").Value = "  -3.38%  "

$ws.Range("D25").Value = "152.73"
$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("D26").Value = "6.62"
$ws.Range("E26").Value = "  -1.44%  "

$ws.Range("D27").Value = "14.82"
$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("E28").Value = "  +0.17%  "

$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("E30").Value = "  -1.59%  "

$ws.Range("E31").Value = "  -0.73%  "

$ws.Range("E32").Value = "  +1.84%  "

$ws.Range("D33").Value = "1.359.09"
$ws.Range("E33").Value = "  -3.20%  "

$ws.Range("E34").Value = "  +0.61%  "

$ws.Range("E35").Value = "  +0.33%  "

$ws.Range("D36").Value = "0.961"
$ws.Range("E36").Value = "  +5.07%  "

$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("E38").Value = "  +0.31%  "

$ws.Range("D39").Value = "0.521"
$ws.Range("E39").Value = "  -1.09%  "

$ws.Range("D40").Value = "0.805"
$ws.Range("E40").Value = "  -0.51%  "

$ws.Range("E41").Value = "  +0.16%  "

$ws.Range("D42").Value = "5.60"

$ws.Range("E43").Value = "  -0.92%  "

$ws.Range("D44").Value = "2.22"
$ws.Range("E44").Value = "  +1.92%  "

$ws.Range("D45").Value = "63.43"
$ws.Range("E45").Value = "  +0.28%  "

$ws.Range("E46").Value = "  -1.97%  "

$ws.Range("D47").Value = "1.680.04"
$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("D48").Value = "86.19"
$ws.Range("E48").Value = "  -0.32%  "

$ws.Range("E49").Value = "  +1.11%  "

$ws.Range("D50").Value = "0.0₇0968"
$ws.Range("E50").Value = "  -1.10%  "

$ws.Range("D51").Value = "0.0948"
$ws.Range("E51").Value = "  +0.04%  "

# Restore the default cell style on column D so formatting matches the source file.
$ws.Range("D2:D51").Style = "Normal"
